$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# Build the new query text (includes a zero-width space character on its own line,
# matching the source content exactly).
$zw = [char]0x200b
$query = "MATCH (ct:clinical_trial)<--(a:arm)<--(c:case)<-[*]-(prt)<--(f:file)`n" +
         "    WHERE c.gender = ""MALE""`n" +
         "WITH DISTINCT f, prt, c, a, ct`n" +
         "RETURN`n" +
         "    COALESCE(f.file_name, '') AS ``File Name``,`n" +
         "    COALESCE(head(labels(prt)), '') AS ``Association``,`n" +
         $zw + "`n" +
         "    COALESCE(f.file_description, '') AS ``Description``,`n" +
         "    COALESCE(f.file_format, '') AS ``File Format``,`n" +
         "    COALESCE(f.file_size, '') AS ``Size``,`n" +
         "    COALESCE(ct.clinical_trial_designation, '') AS ``Trial Code``,`n" +
         "    COALESCE(a.arm_id, '') AS ``Arm``,`n" +
         "    COALESCE(c.case_id, '') AS ``Case ID``"

# Fill in new row 3 values
$ws.Range("A3").Value = "FilesTab"
$ws.Range("B3").Value = $query
$ws.Range("C3").Value = $ws.Range("C2").Value2
$ws.Range("D3").Value = $ws.Range("D2").Value2
$ws.Range("E3").Value = $ws.Range("E2").Value2

# Apply same style (wrap text) as row 2 cells B2/C2 to the new B3/C3 cells
$ws.Range("B3").WrapText = $true
$ws.Range("C3").WrapText = $true

# Set row height for the new row
$ws.Rows.Item(3).RowHeight = 188.5

# Update the active selection to match the edited workbook (B3 selected)
$ws.Range("B3").Select()
